$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new journal entries (rows 15-21) ---
# Values are entered in the same order the original author typed them, so
# that new shared-string entries land at the same indices as the source.

# Row 15
$ws.Range("A15").Value = "19/03/18"
$ws.Range("B15").Value = "Retour sur le travial du weekend"
$ws.Range("C15").Value = 1

# Row 16 (date typed as a real date serial, not text)
$ws.Range("A16").Value = 43192
$ws.Range("B16").Value = "Retour sur le travail des vacances"
$ws.Range("C16").Value = 1

# Row 17 date only for now (also a real date serial)
$ws.Range("A17").Value = 43199

# Row 18
$ws.Range("A18").Value = "15/04/18"
$ws.Range("B18").Value = "travail en groupe chez moi, refactorisation du code, bilan sur l'état du projet, implémentation du crayon"
$ws.Range("C18").Value = 8

# Back to row 17's activity text
$ws.Range("B17").Value = "Discution sur les problèmes encontré, plaification du weekend de travail à venir"

# Row 19
$ws.Range("A19").Value = "16/04/18"
$ws.Range("B19").Value = "retour sur le travail du weekend, présentation de l'état du projet a M.Rentch"
$ws.Range("C19").Value = 1

# Row 20
$ws.Range("A20").Value = "23/04/18"
$ws.Range("B20").Value = "Semaine trop chargé, pas de travail"
$ws.Range("C20").Value = 0

# Row 21
$ws.Range("A21").Value = "29/04/18"
$ws.Range("B21").Value = "Remaniement du modèle MVD, restructuration des dossiers, travail sur le crayon et la gomme et recherches sur la structure ""canevas"" de FXML"
$ws.Range("C21").Value = 9

# Row heights for the wrapped multi-line entries
$ws.Rows.Item(17).RowHeight = 30
$ws.Rows.Item(18).RowHeight = 30
$ws.Rows.Item(19).RowHeight = 30
$ws.Rows.Item(21).RowHeight = 45

# --- Insert a new blank row before the Total row, preserving formatting ---
$ws.Rows.Item(32).Copy()
$ws.Rows.Item(33).Insert(-4121)
$ws.Range("A32:C32").Copy()
$ws.Range("A33:C33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the Total formula to include the new row ---
$ws.Range("C34").Formula = "=SUM(C5:C33)"
